# TST: no person itself in coauthor
#
# Row 37 of the "NSF COA Template" sheet (the data row of the "TableD"
# co-authors table) listed the document owner ("Scopatz, Anthony") as
# their own co-author. Remove that erroneous self-listing by clearing
# the row's contents, leaving the table row blank (same treatment as
# the other, already-empty table rows on this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "C", "D", "E")

foreach ($col in $cols) {
    $cell = $ws.Range($col + "37")

    # Drop the cell's text/number content -- row becomes an empty data row.
    $cell.ClearContents()

    # This was the only filled-in row of the table, so it carried a
    # "closing" border (boxed on all four sides). Once blank it should
    # match the look of the table's other empty rows: no left edge,
    # keep the top/right/bottom rules.
    $cell.Borders.Item(7).LineStyle = 0
}
